$d = $word.ActiveDocument

# Build the replacement OOXML for the whole document body (section properties
# are intentionally omitted; Word preserves the existing sectPr automatically).
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:color w:val="00B050"/>
                <w:sz w:val="56"/>
                <w:szCs w:val="56"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="00B050"/>
                <w:sz w:val="56"/>
                <w:szCs w:val="56"/>
              </w:rPr>
              <w:t>INTRODUCTION</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t>This is the first time I have ever experience online education. I always had strange dreams but I never dreams that there could be a thing like online education. But now that I see, now that I experiencing it, I know that there is a thing named online education. On my first day of experiencing online education, everything was new to me. It was very strange and unusual to me. Now I am used to it. In the holidays, I long for</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t xml:space="preserve"> online school to start, </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t>and  when</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t xml:space="preserve"> the school starts, I want the holidays to start : )</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:b/>
                <w:i/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
                <w:u w:val="single"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t xml:space="preserve">So, </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t>I can&#8217;t really tell which one do I like, school or holidays.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:i/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
                <w:u w:val="single"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:color w:val="7030A0"/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
            </w:pPr>
            <w:hyperlink r:id="rId4" w:history="1">
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                  <w:sz w:val="36"/>
                  <w:szCs w:val="36"/>
                </w:rPr>
                <w:t>GO T</w:t>
              </w:r>
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                  <w:sz w:val="36"/>
                  <w:szCs w:val="36"/>
                </w:rPr>
                <w:t>O</w:t>
              </w:r>
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                  <w:sz w:val="36"/>
                  <w:szCs w:val="36"/>
                </w:rPr>
                <w:t xml:space="preserve"> </w:t>
              </w:r>
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                  <w:sz w:val="36"/>
                  <w:szCs w:val="36"/>
                </w:rPr>
                <w:t>HOME</w:t>
              </w:r>
            </w:hyperlink>
            <w:r>
              <w:rPr>
                <w:rStyle w:val="Hyperlink"/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t xml:space="preserve">  </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:color w:val="FF0000"/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:color w:val="FF0000"/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
                <w:u w:val="single"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rIdStylesHelper" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml">
    <pkg:xmlData>
      <w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:style w:type="character" w:styleId="Hyperlink">
          <w:name w:val="Hyperlink"/>
          <w:rPr>
            <w:color w:val="0563C1"/>
            <w:u w:val="single"/>
          </w:rPr>
        </w:style>
      </w:styles>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml)
